$d = $word.ActiveDocument

# --- Part 1: update version number and citation in the integration paragraph ---
# " v1.x.x (Bashevkin 202x)" -> " v2.0.0 (Bashevkin et al. 202x)"
$d.Content.Find.Execute(
    " v1.x.x (Bashevkin 202", $true, $false, $false, $false, $false,
    $true, 1, $false, " v2.0.0 (Bashevkin et al. 202", 2) | Out-Null

# --- Part 2: append "4. Literature cited" section with two references ---
function Append-Paragraph([string]$text) {
    $r = $word.ActiveDocument.Content
    $r.Collapse(0)
    $r.InsertParagraphAfter()
    $r = $word.ActiveDocument.Content
    $r.Collapse(0)
    if ($text -ne $null -and $text -ne "") {
        $r.InsertAfter($text)
    }
}

Append-Paragraph ""
Append-Paragraph "4. Literature cited"
Append-Paragraph ""
Append-Paragraph "Bashevkin, S. M., J. W. Gaeta, T. X. Nguyen, L. Mitchell, and S. Khana. 2022. LTMRdata: An R package to integrate data from 9 fish surveys in the San Francisco Estuary. v2.2.0. Zenodo. doi:10.5281/zenodo.6048977"
Append-Paragraph ""
Append-Paragraph "Wickham, H., M. Averick, J. Bryan, and others. 2019. Welcome to the Tidyverse. Journal of Open Source Software 4: 1686. doi:10.21105/joss.01686"
